$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.228.05'
$ws.Range("E2").Value = '  -0.42%  '

$ws.Range("D3").Value = '1.656.60'
$ws.Range("E3").Value = '  -0.71%  '

$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  -0.41%  '

$ws.Range("D5").Value = '''219.49'
$ws.Range("E5").Value = '  -0.47%  '

$ws.Range("D6").Value = '''0.5255'
$ws.Range("E6").Value = '  -1.08%  '

$ws.Range("E7").Value = '  -0.36%  '

$ws.Range("D8").Value = '''0.2671'
$ws.Range("E8").Value = '  +0.73%  '

$ws.Range("D9").Value = '''0.06361'
$ws.Range("E9").Value = '  -0.23%  '

$ws.Range("D10").Value = '''20.73'
$ws.Range("E10").Value = '  -0.82%  '

$ws.Range("D11").Value = '''0.07723'
$ws.Range("E11").Value = '  -1.56%  '

$ws.Range("D12").Value = '''4.606'
$ws.Range("E12").Value = '  +1.69%  '

$ws.Range("D13").Value = '1.656.42'
$ws.Range("E13").Value = '  -0.11%  '

$ws.Range("D14").Value = '1.885.34'

$ws.Range("D15").Value = '''0.5641'
$ws.Range("E15").Value = '  +0.63%  '

$ws.Range("D16").Value = '0.0₅8229'
$ws.Range("E16").Value = '  +0.74%  '

$ws.Range("D17").Value = '''65.52'
$ws.Range("E17").Value = '  -0.54%  '

$ws.Range("D18").Value = '26.232.53'
$ws.Range("E18").Value = '  -0.39%  '

$ws.Range("E19").Value = '  -0.38%  '

$ws.Range("D20").Value = '''4.722'
$ws.Range("E20").Value = '  -0.15%  '

$ws.Range("D21").Value = '''10.42'
$ws.Range("E21").Value = '  +1.33%  '

$ws.Range("D22").Value = '''192.24'
$ws.Range("E22").Value = '  -2.83%  '

$ws.Range("D23").Value = '''6.010'
$ws.Range("E23").Value = '  -0.70%  '

$ws.Range("D24").Value = '''1.004'

$ws.Range("D25").Value = '''144.40'
$ws.Range("E25").Value = '  -1.29%  '

$ws.Range("E26").Value = '  -1.12%  '

$ws.Range("D27").Value = '''7.280'
$ws.Range("E27").Value = '  +0.32%  '

$ws.Range("E28").Value = '  -1.32%  '

$ws.Range("D29").Value = '''1.497'
$ws.Range("E29").Value = '  -0.62%  '

$ws.Range("D30").Value = '''0.05637'
$ws.Range("E30").Value = '  -4.53%  '

$ws.Range("D31").Value = '''1.275'
$ws.Range("E31").Value = '  -0.64%  '

$ws.Range("D32").Value = '''3.509'
$ws.Range("E32").Value = '  -1.28%  '

$ws.Range("D33").Value = '''3.380'
$ws.Range("E33").Value = '  +1.54%  '

$ws.Range("D34").Value = '''1.585'
$ws.Range("E34").Value = '  -1.37%  '

$ws.Range("D35").Value = '''0.9548'

$ws.Range("D36").Value = '''2.799'
$ws.Range("E36").Value = '  -1.11%  '

$ws.Range("E37").Value = '  -0.96%  '

$ws.Range("D38").Value = '''0.5769'
$ws.Range("E38").Value = '  -1.07%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.01600'
$ws.Range("E39").Value = '  -1.05%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '''6.015'
$ws.Range("E40").Value = '  +0.98%  '

$ws.Range("E41").Value = '  -0.43%  '

$ws.Range("D42").Value = '''0.8426'
$ws.Range("E42").Value = '  -1.78%  '

$ws.Range("D43").Value = '''101.98'
$ws.Range("E43").Value = '  -0.97%  '

$ws.Range("D44").Value = '1.014.05'
$ws.Range("E44").Value = '  -6.03%  '

$ws.Range("D45").Value = '1.796.03'
$ws.Range("E45").Value = '  -0.61%  '

$ws.Range("D46").Value = '''58.60'
$ws.Range("E46").Value = '  -0.09%  '

$ws.Range("E47").Value = '  -0.87%  '

$ws.Range("D48").Value = '''0.05342'
$ws.Range("E48").Value = '  +3.72%  '

$ws.Range("D49").Value = '''0.4349'
$ws.Range("E49").Value = '  -1.25%  '

$ws.Range("D50").Value = '''8.007'
$ws.Range("E50").Value = '  -0.85%  '

$ws.Range("D51").Value = '''0.09800'
$ws.Range("E51").Value = '  +2.00%  '
